$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038550885001746
$ws.Range("D2").Value = 1.041959221025497
$ws.Range("E2").Value = 1.047458027391036
$ws.Range("F2").Value = 1.059024497997742
$ws.Range("I2").Value = 1.040383642019615
$ws.Range("J2").Value = 1.043647710921769
$ws.Range("K2").Value = 1.044737261897528
$ws.Range("L2").Value = 1.050220609375033
$ws.Range("M2").Value = 1.061755130604798
$ws.Range("N2").Value = 1.045129810379344
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.039411902332398
$ws.Range("D3").Value = 1.042610888768261
$ws.Range("E3").Value = 1.048257099019398
$ws.Range("F3").Value = 1.059982196274657
$ws.Range("I3").Value = 1.040590834075739
$ws.Range("J3").Value = 1.044154011343644
$ws.Range("K3").Value = 1.045200080920532
$ws.Range("L3").Value = 1.050831562315867
$ws.Range("M3").Value = 1.062526613633087
$ws.Range("N3").Value = 1.045636829805891
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.039969622311174
$ws.Range("D4").Value = 1.043033014547638
$ws.Range("E4").Value = 1.048775086412755
$ws.Range("F4").Value = 1.060603117844128
$ws.Range("I4").Value = 1.040724028327045
$ws.Range("J4").Value = 1.044481545763754
$ws.Range("K4").Value = 1.045499322622235
$ws.Range("L4").Value = 1.051227167799729
$ws.Range("M4").Value = 1.063026415651969
$ws.Range("N4").Value = 1.045964829362446
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.040204226179406
$ws.Range("D5").Value = 1.043210583073931
$ws.Range("E5").Value = 1.048993070394696
$ws.Range("F5").Value = 1.060864445246245
$ws.Range("I5").Value = 1.040779813559424
$ws.Range("J5").Value = 1.044619221916631
$ws.Range("K5").Value = 1.045625066836298
$ws.Range("L5").Value = 1.05139354559542
$ws.Range("M5").Value = 1.063236674773808
$ws.Range("N5").Value = 1.046102701031252
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.040243625273822
$ws.Range("D6").Value = 1.043240403828728
$ws.Range("E6").Value = 1.049029683838958
$ws.Range("F6").Value = 1.060908340338797
$ws.Range("I6").Value = 1.040789167833247
$ws.Range("J6").Value = 1.044642337181775
$ws.Range("K6").Value = 1.045646176466507
$ws.Range("L6").Value = 1.051421484947311
$ws.Range("M6").Value = 1.063271986509462
$ws.Range("N6").Value = 1.046125849122724
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.039972756557599
$ws.Range("D7").Value = 1.043035386809198
$ws.Range("E7").Value = 1.048777998254852
$ws.Range("F7").Value = 1.060606608568546
$ws.Range("I7").Value = 1.040724774556585
$ws.Range("J7").Value = 1.044483385476614
$ws.Range("K7").Value = 1.045501003048016
$ws.Range("L7").Value = 1.051229390693245
$ws.Range("M7").Value = 1.06302922458718
$ws.Range("N7").Value = 1.04596667168791
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.038841747995815
$ws.Range("D8").Value = 1.042179360370559
$ws.Range("E8").Value = 1.047727882867548
$ws.Range("F8").Value = 1.059347901579132
$ws.Range("I8").Value = 1.040453843876437
$ws.Range("J8").Value = 1.043818832513397
$ws.Range("K8").Value = 1.044893721319769
$ws.Range("L8").Value = 1.050427025408687
$ws.Range("M8").Value = 1.062015731576987
$ws.Range("N8").Value = 1.045301174983259
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036853305588825
$ws.Range("D9").Value = 1.040674475766044
$ws.Range("E9").Value = 1.045884677198996
$ws.Range("F9").Value = 1.0571393701919
$ws.Range("I9").Value = 1.039969773805326
$ws.Range("J9").Value = 1.042647276385409
$ws.Range("K9").Value = 1.043821879827038
$ws.Range("L9").Value = 1.04901534350123
$ws.Range("M9").Value = 1.060234496357663
$ws.Range("N9").Value = 1.044127955111191
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035530816401335
$ws.Range("D10").Value = 1.039673699568495
$ws.Range("E10").Value = 1.044660835036937
$ws.Range("F10").Value = 1.055673485694443
$ws.Range("I10").Value = 1.039642625287893
$ws.Range("J10").Value = 1.041865952927199
$ws.Range("K10").Value = 1.043106223592916
$ws.Range("L10").Value = 1.048075772207779
$ws.Range("M10").Value = 1.059050235073165
$ws.Range("N10").Value = 1.04334552208405
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034958926781765
$ws.Range("D11").Value = 1.039240961341359
$ws.Range("E11").Value = 1.044132095628839
$ws.Range("F11").Value = 1.055040297482303
$ws.Range("I11").Value = 1.039499922753049
$ws.Range("J11").Value = 1.041527578017944
$ws.Range("K11").Value = 1.042796091550279
$ws.Range("L11").Value = 1.047669311656994
$ws.Range("M11").Value = 1.058538221563895
$ws.Range("N11").Value = 1.043006666643612
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034746616685243
$ws.Range("D12").Value = 1.039080315577475
$ws.Range("E12").Value = 1.043935879092734
$ws.Range("F12").Value = 1.054805337697387
$ws.Range("I12").Value = 1.039446760322833
$ws.Range("J12").Value = 1.041401883159484
$ws.Range("K12").Value = 1.042680858597131
$ws.Range("L12").Value = 1.047518392531552
$ws.Range("M12").Value = 1.058348155563379
$ws.Range("N12").Value = 1.042880793284036
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034792152686536
$ws.Range("D13").Value = 1.039114770427435
$ws.Range("E13").Value = 1.043977959996093
$ws.Range("F13").Value = 1.054855726729089
$ws.Range("I13").Value = 1.039458170909433
$ws.Range("J13").Value = 1.041428845447205
$ws.Range("K13").Value = 1.042705578071556
$ws.Range("L13").Value = 1.047550762531281
$ws.Range("M13").Value = 1.058388919973837
$ws.Range("N13").Value = 1.042907793861298
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034941374805077
$ws.Range("D14").Value = 1.039227680416179
$ws.Range("E14").Value = 1.044115872610268
$ws.Range("F14").Value = 1.055020870839283
$ws.Range("I14").Value = 1.039495531515782
$ws.Range("J14").Value = 1.041517188184938
$ws.Range("K14").Value = 1.0427865670881
$ws.Range("L14").Value = 1.047656835435559
$ws.Range("M14").Value = 1.058522508220552
$ws.Range("N14").Value = 1.042996262055851
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035033330782124
$ws.Range("D15").Value = 1.039297260291242
$ws.Range("E15").Value = 1.044200869134774
$ws.Range("F15").Value = 1.055122652702422
$ws.Range("I15").Value = 1.039518529919486
$ws.Range("J15").Value = 1.041571618122179
$ws.Range("K15").Value = 1.042836462345392
$ws.Range("L15").Value = 1.047722198233774
$ws.Range("M15").Value = 1.05860483211113
$ws.Range("N15").Value = 1.043050769289847
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.03556878690331
$ws.Range("D16").Value = 1.039702431857219
$ws.Range("E16").Value = 1.044695950995146
$ws.Range("F16").Value = 1.055715541102799
$ws.Range("I16").Value = 1.0396520740226
$ws.Range("J16").Value = 1.041888408660572
$ws.Range("K16").Value = 1.043126800907305
$ws.Range("L16").Value = 1.048102755775715
$ws.Range("M16").Value = 1.059084232274436
$ws.Range("N16").Value = 1.04336800970714
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035904868055733
$ws.Range("D17").Value = 1.039956748286822
$ws.Range("E17").Value = 1.045006823177405
$ws.Range("F17").Value = 1.05608786037891
$ws.Range("I17").Value = 1.039735563377356
$ws.Range("J17").Value = 1.042087108640041
$ws.Range("K17").Value = 1.043308857007616
$ws.Range("L17").Value = 1.048341572146026
$ws.Range("M17").Value = 1.059385157076463
$ws.Range("N17").Value = 1.043566991863373
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.036100971440642
$ws.Range("D18").Value = 1.04010514496847
$ws.Range("E18").Value = 1.045188264590474
$ws.Range("F18").Value = 1.056305177297146
$ws.Range("I18").Value = 1.03978416042752
$ws.Range("J18").Value = 1.042203001427669
$ws.Range("K18").Value = 1.043415023226667
$ws.Range("L18").Value = 1.04848090627167
$ws.Range("M18").Value = 1.059560756525535
$ws.Range("N18").Value = 1.043683049232051
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.036167849954015
$ws.Range("D19").Value = 1.040155754224866
$ws.Range("E19").Value = 1.04525015086599
$ws.Range("F19").Value = 1.056379302035986
$ws.Range("I19").Value = 1.039800713635233
$ws.Range("J19").Value = 1.042242516896682
$ws.Range("K19").Value = 1.04345121905296
$ws.Range("L19").Value = 1.048528421771586
$ws.Range("M19").Value = 1.059620644092637
$ws.Range("N19").Value = 1.043722620817563
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035868802174251
$ws.Range("D20").Value = 1.039929456507693
$ws.Range("E20").Value = 1.044973457634509
$ws.Range("F20").Value = 1.056047898580294
$ws.Range("I20").Value = 1.039726616189282
$ws.Range("J20").Value = 1.042065790578222
$ws.Range("K20").Value = 1.043289326587417
$ws.Range("L20").Value = 1.048315945616687
$ws.Range("M20").Value = 1.05935286291393
$ws.Range("N20").Value = 1.04354564352746
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.034897429419448
$ws.Range("D21").Value = 1.039194428672386
$ws.Range("E21").Value = 1.04407525577763
$ws.Range("F21").Value = 1.054972233527472
$ws.Range("I21").Value = 1.039484534059293
$ws.Range("J21").Value = 1.041491173644107
$ws.Range("K21").Value = 1.042762718828481
$ws.Range("L21").Value = 1.047625597985683
$ws.Range("M21").Value = 1.058483166525006
$ws.Range("N21").Value = 1.042970210571389
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034287355884106
$ws.Range("D22").Value = 1.038732823391668
$ws.Range("E22").Value = 1.0435115676386
$ws.Range("F22").Value = 1.054297277872092
$ws.Range("I22").Value = 1.039331423120282
$ws.Range("J22").Value = 1.041129847219045
$ws.Range("K22").Value = 1.042431411328752
$ws.Range("L22").Value = 1.047191887692833
$ws.Range("M22").Value = 1.057937040902904
$ws.Range("N22").Value = 1.042608371021359
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.034610703691659
$ws.Range("D23").Value = 1.038977477732299
$ws.Range("E23").Value = 1.043810289577095
$ws.Range("F23").Value = 1.054654955369389
$ws.Range("I23").Value = 1.039412675623546
$ws.Range("J23").Value = 1.041321396780063
$ws.Range("K23").Value = 1.042607063052827
$ws.Range("L23").Value = 1.047421773180317
$ws.Range("M23").Value = 1.058226486801495
$ws.Range("N23").Value = 1.042800192604724
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03588509856846
$ws.Range("D24").Value = 1.039941788308268
$ws.Range("E24").Value = 1.044988533732238
$ws.Range("F24").Value = 1.056065955135585
$ws.Range("I24").Value = 1.039730659350218
$ws.Range("J24").Value = 1.042075423310015
$ws.Range("K24").Value = 1.043298151618
$ws.Range("L24").Value = 1.048327525029354
$ws.Range("M24").Value = 1.059367455024012
$ws.Range("N24").Value = 1.043555289938837
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03736681872841
$ws.Range("D25").Value = 1.041063094515514
$ws.Range("E25").Value = 1.046360322803875
$ws.Range("F25").Value = 1.057709195626519
$ws.Range("I25").Value = 1.040095702161124
$ws.Range("J25").Value = 1.042950207240976
$ws.Range("K25").Value = 1.044099174736212
$ws.Range("L25").Value = 1.049380029687406
$ws.Range("M25").Value = 1.060694425827881
$ws.Range("N25").Value = 1.044431316163319
